# Applies the "Otomatik guncelleme: 2025-06-21 02:29:20" data refresh.
#
# Summary of the change:
#  - "eskalasyon": per Name-group rolling window of observations shifts down
#    by one row (oldest date dropped for each group) and a new 2025-06-21
#    (serial 45829) observation is inserted at the top of every group; the
#    table grows from 35 to 44 data rows (A1:G36 -> A1:G45).
#  - "durum": StartDate column bumped from the literal text "2025-06-19" to
#    "2025-06-21" for every district row.
#  - District price-history sheets ("934015","065001","035001","055001",
#    "021001","038001"): one new trailing row appended for priceDate
#    2025-06-21 (serial 45829) with the refreshed "amount".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "eskalasyon" sheet - rewrite the full A2:G45 data block
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("eskalasyon")
$dateFmt = $ws.Cells.Item(2,1).NumberFormat

$ws.Cells.Item(2,1).Value = 45673
$ws.Cells.Item(2,1).NumberFormat = $dateFmt
$ws.Cells.Item(2,2).Value = 'Motorin UltraForce'
$ws.Cells.Item(2,3).Value = 49.23
$ws.Cells.Item(2,4).Value = 0.06374243733794294
$ws.Cells.Item(2,5).Value = 0.06374243733794294
$ws.Cells.Item(2,6).Value = 'Servis Diyarbakır'
$ws.Cells.Item(2,7).Value = 0.05
$ws.Cells.Item(3,1).Value = 45756
$ws.Cells.Item(3,1).NumberFormat = $dateFmt
$ws.Cells.Item(3,2).Value = 'Motorin UltraForce'
$ws.Cells.Item(3,3).Value = 46.38
$ws.Cells.Item(3,4).Value = -0.0578915295551492
$ws.Cells.Item(3,5).Value = -0.0578915295551492
$ws.Cells.Item(3,6).Value = 'Servis Diyarbakır'
$ws.Cells.Item(3,7).Value = 0.05
$ws.Cells.Item(4,1).Value = 45822
$ws.Cells.Item(4,1).NumberFormat = $dateFmt
$ws.Cells.Item(4,2).Value = 'Motorin UltraForce'
$ws.Cells.Item(4,3).Value = 49.96
$ws.Cells.Item(4,4).Value = 0.07718844329452357
$ws.Cells.Item(4,5).Value = 0.07718844329452357
$ws.Cells.Item(4,6).Value = 'Servis Diyarbakır'
$ws.Cells.Item(4,7).Value = 0.05
$ws.Cells.Item(5,1).Value = 45827
$ws.Cells.Item(5,1).NumberFormat = $dateFmt
$ws.Cells.Item(5,2).Value = 'Motorin UltraForce'
$ws.Cells.Item(5,3).Value = 53.57
$ws.Cells.Item(5,4).Value = 0.07225780624499589
$ws.Cells.Item(5,5).Value = 0.07225780624499589
$ws.Cells.Item(5,6).Value = 'Servis Diyarbakır'
$ws.Cells.Item(5,7).Value = 0.05
$ws.Cells.Item(6,1).Value = 45829
$ws.Cells.Item(6,1).NumberFormat = $dateFmt
$ws.Cells.Item(6,2).Value = 'Motorin UltraForce'
$ws.Cells.Item(6,3).Value = 56.91
$ws.Cells.Item(6,4).Value = 0.06234832928878098
$ws.Cells.Item(6,5).Value = 0.06234832928878098
$ws.Cells.Item(6,6).Value = 'Servis Diyarbakır'
$ws.Cells.Item(6,7).Value = 0.05
$ws.Cells.Item(7,1).Value = 45673
$ws.Cells.Item(7,1).NumberFormat = $dateFmt
$ws.Cells.Item(7,2).Value = 'Motorin UltraForce'
$ws.Cells.Item(7,3).Value = 48.88
$ws.Cells.Item(7,4).Value = 0.05663640293990491
$ws.Cells.Item(7,5).Value = 0.05663640293990491
$ws.Cells.Item(7,6).Value = 'Servis Kayseri'
$ws.Cells.Item(7,7).Value = 0.05
$ws.Cells.Item(8,1).Value = 45756
$ws.Cells.Item(8,1).NumberFormat = $dateFmt
$ws.Cells.Item(8,2).Value = 'Motorin UltraForce'
$ws.Cells.Item(8,3).Value = 46.03
$ws.Cells.Item(8,4).Value = -0.05830605564648117
$ws.Cells.Item(8,5).Value = -0.05830605564648117
$ws.Cells.Item(8,6).Value = 'Servis Kayseri'
$ws.Cells.Item(8,7).Value = 0.05
$ws.Cells.Item(9,1).Value = 45822
$ws.Cells.Item(9,1).NumberFormat = $dateFmt
$ws.Cells.Item(9,2).Value = 'Motorin UltraForce'
$ws.Cells.Item(9,3).Value = 49.64
$ws.Cells.Item(9,4).Value = 0.07842711275255265
$ws.Cells.Item(9,5).Value = 0.07842711275255265
$ws.Cells.Item(9,6).Value = 'Servis Kayseri'
$ws.Cells.Item(9,7).Value = 0.05
$ws.Cells.Item(10,1).Value = 45827
$ws.Cells.Item(10,1).NumberFormat = $dateFmt
$ws.Cells.Item(10,2).Value = 'Motorin UltraForce'
$ws.Cells.Item(10,3).Value = 53.25
$ws.Cells.Item(10,4).Value = 0.07272360999194194
$ws.Cells.Item(10,5).Value = 0.07272360999194194
$ws.Cells.Item(10,6).Value = 'Servis Kayseri'
$ws.Cells.Item(10,7).Value = 0.05
$ws.Cells.Item(11,1).Value = 45829
$ws.Cells.Item(11,1).NumberFormat = $dateFmt
$ws.Cells.Item(11,2).Value = 'Motorin UltraForce'
$ws.Cells.Item(11,3).Value = 56.59
$ws.Cells.Item(11,4).Value = 0.06272300469483572
$ws.Cells.Item(11,5).Value = 0.06272300469483572
$ws.Cells.Item(11,6).Value = 'Servis Kayseri'
$ws.Cells.Item(11,7).Value = 0.05
$ws.Cells.Item(12,1).Value = 45673
$ws.Cells.Item(12,1).NumberFormat = $dateFmt
$ws.Cells.Item(12,2).Value = 'Motorin UltraForce'
$ws.Cells.Item(12,3).Value = 48.57
$ws.Cells.Item(12,4).Value = 0.05678851174934718
$ws.Cells.Item(12,5).Value = 0.05678851174934718
$ws.Cells.Item(12,6).Value = 'Servis Samsun'
$ws.Cells.Item(12,7).Value = 0.05
$ws.Cells.Item(13,1).Value = 45756
$ws.Cells.Item(13,1).NumberFormat = $dateFmt
$ws.Cells.Item(13,2).Value = 'Motorin UltraForce'
$ws.Cells.Item(13,3).Value = 45.66
$ws.Cells.Item(13,4).Value = -0.05991352686843743
$ws.Cells.Item(13,5).Value = -0.05991352686843743
$ws.Cells.Item(13,6).Value = 'Servis Samsun'
$ws.Cells.Item(13,7).Value = 0.05
$ws.Cells.Item(14,1).Value = 45822
$ws.Cells.Item(14,1).NumberFormat = $dateFmt
$ws.Cells.Item(14,2).Value = 'Motorin UltraForce'
$ws.Cells.Item(14,3).Value = 49.24
$ws.Cells.Item(14,4).Value = 0.07840560665790641
$ws.Cells.Item(14,5).Value = 0.07840560665790641
$ws.Cells.Item(14,6).Value = 'Servis Samsun'
$ws.Cells.Item(14,7).Value = 0.05
$ws.Cells.Item(15,1).Value = 45827
$ws.Cells.Item(15,1).NumberFormat = $dateFmt
$ws.Cells.Item(15,2).Value = 'Motorin UltraForce'
$ws.Cells.Item(15,3).Value = 52.85
$ws.Cells.Item(15,4).Value = 0.07331437855402112
$ws.Cells.Item(15,5).Value = 0.07331437855402112
$ws.Cells.Item(15,6).Value = 'Servis Samsun'
$ws.Cells.Item(15,7).Value = 0.05
$ws.Cells.Item(16,1).Value = 45829
$ws.Cells.Item(16,1).NumberFormat = $dateFmt
$ws.Cells.Item(16,2).Value = 'Motorin UltraForce'
$ws.Cells.Item(16,3).Value = 56.19
$ws.Cells.Item(16,4).Value = 0.06319772942289492
$ws.Cells.Item(16,5).Value = 0.06319772942289492
$ws.Cells.Item(16,6).Value = 'Servis Samsun'
$ws.Cells.Item(16,7).Value = 0.05
$ws.Cells.Item(17,1).Value = 45673
$ws.Cells.Item(17,1).NumberFormat = $dateFmt
$ws.Cells.Item(17,2).Value = 'Motorin UltraForce'
$ws.Cells.Item(17,3).Value = 47.24
$ws.Cells.Item(17,4).Value = 0.07023108291798819
$ws.Cells.Item(17,5).Value = 0.07023108291798819
$ws.Cells.Item(17,6).Value = 'Spot Araç Anadolu Toplama'
$ws.Cells.Item(17,7).Value = 0.05
$ws.Cells.Item(18,1).Value = 45756
$ws.Cells.Item(18,1).NumberFormat = $dateFmt
$ws.Cells.Item(18,2).Value = 'Motorin UltraForce'
$ws.Cells.Item(18,3).Value = 44.26
$ws.Cells.Item(18,4).Value = -0.06308213378492811
$ws.Cells.Item(18,5).Value = -0.06308213378492811
$ws.Cells.Item(18,6).Value = 'Spot Araç Anadolu Toplama'
$ws.Cells.Item(18,7).Value = 0.05
$ws.Cells.Item(19,1).Value = 45822
$ws.Cells.Item(19,1).NumberFormat = $dateFmt
$ws.Cells.Item(19,2).Value = 'Motorin UltraForce'
$ws.Cells.Item(19,3).Value = 47.76
$ws.Cells.Item(19,4).Value = 0.07907817442385912
$ws.Cells.Item(19,5).Value = 0.07907817442385912
$ws.Cells.Item(19,6).Value = 'Spot Araç Anadolu Toplama'
$ws.Cells.Item(19,7).Value = 0.05
$ws.Cells.Item(20,1).Value = 45827
$ws.Cells.Item(20,1).NumberFormat = $dateFmt
$ws.Cells.Item(20,2).Value = 'Motorin UltraForce'
$ws.Cells.Item(20,3).Value = 51.37
$ws.Cells.Item(20,4).Value = 0.07558626465661633
$ws.Cells.Item(20,5).Value = 0.07558626465661633
$ws.Cells.Item(20,6).Value = 'Spot Araç Anadolu Toplama'
$ws.Cells.Item(20,7).Value = 0.05
$ws.Cells.Item(21,1).Value = 45829
$ws.Cells.Item(21,1).NumberFormat = $dateFmt
$ws.Cells.Item(21,2).Value = 'Motorin UltraForce'
$ws.Cells.Item(21,3).Value = 54.71
$ws.Cells.Item(21,4).Value = 0.06501849328401788
$ws.Cells.Item(21,5).Value = 0.06501849328401788
$ws.Cells.Item(21,6).Value = 'Spot Araç Anadolu Toplama'
$ws.Cells.Item(21,7).Value = 0.05
$ws.Cells.Item(22,1).Value = 45673
$ws.Cells.Item(22,1).NumberFormat = $dateFmt
$ws.Cells.Item(22,2).Value = 'Motorin UltraForce'
$ws.Cells.Item(22,3).Value = 47.24
$ws.Cells.Item(22,4).Value = 0.07023108291798819
$ws.Cells.Item(22,5).Value = 0.07023108291798819
$ws.Cells.Item(22,6).Value = 'Spot Araç Avrupa&Anadolu'
$ws.Cells.Item(22,7).Value = 0.05
$ws.Cells.Item(23,1).Value = 45756
$ws.Cells.Item(23,1).NumberFormat = $dateFmt
$ws.Cells.Item(23,2).Value = 'Motorin UltraForce'
$ws.Cells.Item(23,3).Value = 44.26
$ws.Cells.Item(23,4).Value = -0.06308213378492811
$ws.Cells.Item(23,5).Value = -0.06308213378492811
$ws.Cells.Item(23,6).Value = 'Spot Araç Avrupa&Anadolu'
$ws.Cells.Item(23,7).Value = 0.05
$ws.Cells.Item(24,1).Value = 45822
$ws.Cells.Item(24,1).NumberFormat = $dateFmt
$ws.Cells.Item(24,2).Value = 'Motorin UltraForce'
$ws.Cells.Item(24,3).Value = 47.76
$ws.Cells.Item(24,4).Value = 0.07907817442385912
$ws.Cells.Item(24,5).Value = 0.07907817442385912
$ws.Cells.Item(24,6).Value = 'Spot Araç Avrupa&Anadolu'
$ws.Cells.Item(24,7).Value = 0.05
$ws.Cells.Item(25,1).Value = 45827
$ws.Cells.Item(25,1).NumberFormat = $dateFmt
$ws.Cells.Item(25,2).Value = 'Motorin UltraForce'
$ws.Cells.Item(25,3).Value = 51.37
$ws.Cells.Item(25,4).Value = 0.07558626465661633
$ws.Cells.Item(25,5).Value = 0.07558626465661633
$ws.Cells.Item(25,6).Value = 'Spot Araç Avrupa&Anadolu'
$ws.Cells.Item(25,7).Value = 0.05
$ws.Cells.Item(26,1).Value = 45829
$ws.Cells.Item(26,1).NumberFormat = $dateFmt
$ws.Cells.Item(26,2).Value = 'Motorin UltraForce'
$ws.Cells.Item(26,3).Value = 54.71
$ws.Cells.Item(26,4).Value = 0.06501849328401788
$ws.Cells.Item(26,5).Value = 0.06501849328401788
$ws.Cells.Item(26,6).Value = 'Spot Araç Avrupa&Anadolu'
$ws.Cells.Item(26,7).Value = 0.05
$ws.Cells.Item(27,1).Value = 45673
$ws.Cells.Item(27,1).NumberFormat = $dateFmt
$ws.Cells.Item(27,2).Value = 'Motorin UltraForce'
$ws.Cells.Item(27,3).Value = 47.24
$ws.Cells.Item(27,4).Value = 0.07023108291798819
$ws.Cells.Item(27,5).Value = 0.07023108291798819
$ws.Cells.Item(27,6).Value = 'Spot Araç Teknosa'
$ws.Cells.Item(27,7).Value = 0.05
$ws.Cells.Item(28,1).Value = 45756
$ws.Cells.Item(28,1).NumberFormat = $dateFmt
$ws.Cells.Item(28,2).Value = 'Motorin UltraForce'
$ws.Cells.Item(28,3).Value = 44.26
$ws.Cells.Item(28,4).Value = -0.06308213378492811
$ws.Cells.Item(28,5).Value = -0.06308213378492811
$ws.Cells.Item(28,6).Value = 'Spot Araç Teknosa'
$ws.Cells.Item(28,7).Value = 0.05
$ws.Cells.Item(29,1).Value = 45822
$ws.Cells.Item(29,1).NumberFormat = $dateFmt
$ws.Cells.Item(29,2).Value = 'Motorin UltraForce'
$ws.Cells.Item(29,3).Value = 47.76
$ws.Cells.Item(29,4).Value = 0.07907817442385912
$ws.Cells.Item(29,5).Value = 0.07907817442385912
$ws.Cells.Item(29,6).Value = 'Spot Araç Teknosa'
$ws.Cells.Item(29,7).Value = 0.05
$ws.Cells.Item(30,1).Value = 45827
$ws.Cells.Item(30,1).NumberFormat = $dateFmt
$ws.Cells.Item(30,2).Value = 'Motorin UltraForce'
$ws.Cells.Item(30,3).Value = 51.37
$ws.Cells.Item(30,4).Value = 0.07558626465661633
$ws.Cells.Item(30,5).Value = 0.07558626465661633
$ws.Cells.Item(30,6).Value = 'Spot Araç Teknosa'
$ws.Cells.Item(30,7).Value = 0.05
$ws.Cells.Item(31,1).Value = 45829
$ws.Cells.Item(31,1).NumberFormat = $dateFmt
$ws.Cells.Item(31,2).Value = 'Motorin UltraForce'
$ws.Cells.Item(31,3).Value = 54.71
$ws.Cells.Item(31,4).Value = 0.06501849328401788
$ws.Cells.Item(31,5).Value = 0.06501849328401788
$ws.Cells.Item(31,6).Value = 'Spot Araç Teknosa'
$ws.Cells.Item(31,7).Value = 0.05
$ws.Cells.Item(32,1).Value = 45673
$ws.Cells.Item(32,1).NumberFormat = $dateFmt
$ws.Cells.Item(32,2).Value = 'Motorin UltraForce'
$ws.Cells.Item(32,3).Value = 47.24
$ws.Cells.Item(32,4).Value = 0.07023108291798819
$ws.Cells.Item(32,5).Value = 0.07023108291798819
$ws.Cells.Item(32,6).Value = 'TL/Desi Avrupa Toplama'
$ws.Cells.Item(32,7).Value = 0.05
$ws.Cells.Item(33,1).Value = 45756
$ws.Cells.Item(33,1).NumberFormat = $dateFmt
$ws.Cells.Item(33,2).Value = 'Motorin UltraForce'
$ws.Cells.Item(33,3).Value = 44.26
$ws.Cells.Item(33,4).Value = -0.06308213378492811
$ws.Cells.Item(33,5).Value = -0.06308213378492811
$ws.Cells.Item(33,6).Value = 'TL/Desi Avrupa Toplama'
$ws.Cells.Item(33,7).Value = 0.05
$ws.Cells.Item(34,1).Value = 45822
$ws.Cells.Item(34,1).NumberFormat = $dateFmt
$ws.Cells.Item(34,2).Value = 'Motorin UltraForce'
$ws.Cells.Item(34,3).Value = 47.76
$ws.Cells.Item(34,4).Value = 0.07907817442385912
$ws.Cells.Item(34,5).Value = 0.07907817442385912
$ws.Cells.Item(34,6).Value = 'TL/Desi Avrupa Toplama'
$ws.Cells.Item(34,7).Value = 0.05
$ws.Cells.Item(35,1).Value = 45827
$ws.Cells.Item(35,1).NumberFormat = $dateFmt
$ws.Cells.Item(35,2).Value = 'Motorin UltraForce'
$ws.Cells.Item(35,3).Value = 51.37
$ws.Cells.Item(35,4).Value = 0.07558626465661633
$ws.Cells.Item(35,5).Value = 0.07558626465661633
$ws.Cells.Item(35,6).Value = 'TL/Desi Avrupa Toplama'
$ws.Cells.Item(35,7).Value = 0.05
$ws.Cells.Item(36,1).Value = 45829
$ws.Cells.Item(36,1).NumberFormat = $dateFmt
$ws.Cells.Item(36,2).Value = 'Motorin UltraForce'
$ws.Cells.Item(36,3).Value = 54.71
$ws.Cells.Item(36,4).Value = 0.06501849328401788
$ws.Cells.Item(36,5).Value = 0.06501849328401788
$ws.Cells.Item(36,6).Value = 'TL/Desi Avrupa Toplama'
$ws.Cells.Item(36,7).Value = 0.05
$ws.Cells.Item(37,1).Value = 45784
$ws.Cells.Item(37,1).NumberFormat = $dateFmt
$ws.Cells.Item(37,2).Value = 'Motorin UltraForce'
$ws.Cells.Item(37,3).Value = 44.2
$ws.Cells.Item(37,4).Value = -0.05089113162980452
$ws.Cells.Item(37,5).Value = -0.05089113162980452
$ws.Cells.Item(37,6).Value = 'TL/Desi Avrupa İade Toplama'
$ws.Cells.Item(37,7).Value = 0.05
$ws.Cells.Item(38,1).Value = 45822
$ws.Cells.Item(38,1).NumberFormat = $dateFmt
$ws.Cells.Item(38,2).Value = 'Motorin UltraForce'
$ws.Cells.Item(38,3).Value = 47.76
$ws.Cells.Item(38,4).Value = 0.0805429864253393
$ws.Cells.Item(38,5).Value = 0.0805429864253393
$ws.Cells.Item(38,6).Value = 'TL/Desi Avrupa İade Toplama'
$ws.Cells.Item(38,7).Value = 0.05
$ws.Cells.Item(39,1).Value = 45827
$ws.Cells.Item(39,1).NumberFormat = $dateFmt
$ws.Cells.Item(39,2).Value = 'Motorin UltraForce'
$ws.Cells.Item(39,3).Value = 51.37
$ws.Cells.Item(39,4).Value = 0.07558626465661633
$ws.Cells.Item(39,5).Value = 0.07558626465661633
$ws.Cells.Item(39,6).Value = 'TL/Desi Avrupa İade Toplama'
$ws.Cells.Item(39,7).Value = 0.05
$ws.Cells.Item(40,1).Value = 45829
$ws.Cells.Item(40,1).NumberFormat = $dateFmt
$ws.Cells.Item(40,2).Value = 'Motorin UltraForce'
$ws.Cells.Item(40,3).Value = 54.71
$ws.Cells.Item(40,4).Value = 0.06501849328401788
$ws.Cells.Item(40,5).Value = 0.06501849328401788
$ws.Cells.Item(40,6).Value = 'TL/Desi Avrupa İade Toplama'
$ws.Cells.Item(40,7).Value = 0.05
$ws.Cells.Item(41,1).Value = 45673
$ws.Cells.Item(41,1).NumberFormat = $dateFmt
$ws.Cells.Item(41,2).Value = 'Motorin UltraForce'
$ws.Cells.Item(41,3).Value = 47.24
$ws.Cells.Item(41,4).Value = 0.05587840858292359
$ws.Cells.Item(41,5).Value = 0.05587840858292359
$ws.Cells.Item(41,6).Value = 'TL/Desi Avrupa&Anadolu Dağıtım'
$ws.Cells.Item(41,7).Value = 0.05
$ws.Cells.Item(42,1).Value = 45756
$ws.Cells.Item(42,1).NumberFormat = $dateFmt
$ws.Cells.Item(42,2).Value = 'Motorin UltraForce'
$ws.Cells.Item(42,3).Value = 44.26
$ws.Cells.Item(42,4).Value = -0.06308213378492811
$ws.Cells.Item(42,5).Value = -0.06308213378492811
$ws.Cells.Item(42,6).Value = 'TL/Desi Avrupa&Anadolu Dağıtım'
$ws.Cells.Item(42,7).Value = 0.05
$ws.Cells.Item(43,1).Value = 45822
$ws.Cells.Item(43,1).NumberFormat = $dateFmt
$ws.Cells.Item(43,2).Value = 'Motorin UltraForce'
$ws.Cells.Item(43,3).Value = 47.76
$ws.Cells.Item(43,4).Value = 0.07907817442385912
$ws.Cells.Item(43,5).Value = 0.07907817442385912
$ws.Cells.Item(43,6).Value = 'TL/Desi Avrupa&Anadolu Dağıtım'
$ws.Cells.Item(43,7).Value = 0.05
$ws.Cells.Item(44,1).Value = 45827
$ws.Cells.Item(44,1).NumberFormat = $dateFmt
$ws.Cells.Item(44,2).Value = 'Motorin UltraForce'
$ws.Cells.Item(44,3).Value = 51.37
$ws.Cells.Item(44,4).Value = 0.07558626465661633
$ws.Cells.Item(44,5).Value = 0.07558626465661633
$ws.Cells.Item(44,6).Value = 'TL/Desi Avrupa&Anadolu Dağıtım'
$ws.Cells.Item(44,7).Value = 0.05
$ws.Cells.Item(45,1).Value = 45829
$ws.Cells.Item(45,1).NumberFormat = $dateFmt
$ws.Cells.Item(45,2).Value = 'Motorin UltraForce'
$ws.Cells.Item(45,3).Value = 54.71
$ws.Cells.Item(45,4).Value = 0.06501849328401788
$ws.Cells.Item(45,5).Value = 0.06501849328401788
$ws.Cells.Item(45,6).Value = 'TL/Desi Avrupa&Anadolu Dağıtım'
$ws.Cells.Item(45,7).Value = 0.05

# ---------------------------------------------------------------------
# 2) "durum" sheet - bump the StartDate text column
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("durum")
$ws.Cells.Item(2,2).Value = "'2025-06-21"
$ws.Cells.Item(3,2).Value = "'2025-06-21"
$ws.Cells.Item(4,2).Value = "'2025-06-21"
$ws.Cells.Item(5,2).Value = "'2025-06-21"
$ws.Cells.Item(6,2).Value = "'2025-06-21"
$ws.Cells.Item(7,2).Value = "'2025-06-21"

# ---------------------------------------------------------------------
# 3) District price-history sheets - append the 2025-06-21 observation
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("934015")
$ws.Cells.Item(196,1).Value = 45829
$ws.Cells.Item(196,1).NumberFormat = $ws.Cells.Item(195,1).NumberFormat
$ws.Cells.Item(196,2).Value = 'Motorin UltraForce'
$ws.Cells.Item(196,3).Value = 54.71

$ws = $wb.Worksheets.Item("065001")
$ws.Cells.Item(196,1).Value = 45829
$ws.Cells.Item(196,1).NumberFormat = $ws.Cells.Item(195,1).NumberFormat
$ws.Cells.Item(196,2).Value = 'Motorin UltraForce'
$ws.Cells.Item(196,3).Value = 56.96

$ws = $wb.Worksheets.Item("035001")
$ws.Cells.Item(196,1).Value = 45829
$ws.Cells.Item(196,1).NumberFormat = $ws.Cells.Item(195,1).NumberFormat
$ws.Cells.Item(196,2).Value = 'Motorin UltraForce'
$ws.Cells.Item(196,3).Value = 55.91

$ws = $wb.Worksheets.Item("055001")
$ws.Cells.Item(196,1).Value = 45829
$ws.Cells.Item(196,1).NumberFormat = $ws.Cells.Item(195,1).NumberFormat
$ws.Cells.Item(196,2).Value = 'Motorin UltraForce'
$ws.Cells.Item(196,3).Value = 56.19

$ws = $wb.Worksheets.Item("021001")
$ws.Cells.Item(196,1).Value = 45829
$ws.Cells.Item(196,1).NumberFormat = $ws.Cells.Item(195,1).NumberFormat
$ws.Cells.Item(196,2).Value = 'Motorin UltraForce'
$ws.Cells.Item(196,3).Value = 56.91

$ws = $wb.Worksheets.Item("038001")
$ws.Cells.Item(196,1).Value = 45829
$ws.Cells.Item(196,1).NumberFormat = $ws.Cells.Item(195,1).NumberFormat
$ws.Cells.Item(196,2).Value = 'Motorin UltraForce'
$ws.Cells.Item(196,3).Value = 56.59

